$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TS")

# New attribute label used for the DH solar heating constraint per time-slice
$ws.Range("B22").Value = "R"
$ws.Range("C22").Value = "AF"
$ws.Range("D22").Value = 0.10829415070669778
$ws.Range("E22").Value = "EHSOLHT1E"

$ws.Range("B23").Value = "S"
$ws.Range("C23").Value = "AF"
$ws.Range("D23").Value = 0.12417253524194671
$ws.Range("E23").Value = "EHSOLHT1E"

$ws.Range("B24").Value = "F"
$ws.Range("C24").Value = "AF"
$ws.Range("D24").Value = 0.042282889155662889
$ws.Range("E24").Value = "EHSOLHT1E"

$ws.Range("B25").Value = "W"
$ws.Range("C25").Value = "AF"
$ws.Range("D25").Value = 0.030775520587982483
$ws.Range("E25").Value = "EHSOLHT1E"

# Match formatting of the equivalent existing attribute block (rows 10-13)
$ws.Range("D22:D25").NumberFormat = "0%"
$ws.Range("D22:D25").Interior.Color = 5296274

$ws.Range("I25").Select()
